$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 823.0448
$ws.Range("J17").Value = 824.50793
$ws.Range("L17").Value = 2473.52379
$ws.Range("N17").Value = -2809.52379

$ws.Range("H33").Value = 4400521
$ws.Range("I33").Value = 8250246
$ws.Range("K33").Value = 8250246
$ws.Range("M33").Value = -8250017

$ws.Range("H43").Value = 3434.2307
$ws.Range("J43").Value = 3806.7144
$ws.Range("L43").Value = 3806.7144
$ws.Range("N43").Value = -3944.7144

$ws.Range("H47").Value = 15037
$ws.Range("J47").Value = 20074
$ws.Range("L47").Value = 20074
$ws.Range("N47").Value = -22018

$ws.Range("H61").Value = 387.8
$ws.Range("I61").Value = 387.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1163.4
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -991.4000000000001
$ws.Range("N61").ClearContents()

$ws.Range("H62").Value = 4898.375
$ws.Range("J62").Value = 4999
$ws.Range("L62").Value = 4999
$ws.Range("N62").Value = -6247

$ws.Range("H64").Value = 4738.579
$ws.Range("J64").Value = 5430
$ws.Range("L64").Value = 5430
$ws.Range("N64").Value = -5926

$ws.Range("H65").Value = 4898.375
$ws.Range("J65").Value = 4999
$ws.Range("L65").Value = 24995
$ws.Range("N65").Value = -31235

$ws.Range("H67").Value = 4738.579
$ws.Range("J67").Value = 5430
$ws.Range("L67").Value = 5430
$ws.Range("N67").Value = -7146

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H94").Value = 802.6667

$ws.Range("H107").Value = 437.8889
$ws.Range("I107").Value = 484.66666
$ws.Range("K107").Value = 484.66666
$ws.Range("M107").Value = 1435.33334

$ws.Range("H116").Value = 58302.6
$ws.Range("I116").Value = 8410.875
$ws.Range("K116").Value = 8410.875
$ws.Range("M116").Value = -4968.875

$ws.Range("H125").Value = 849.6667
$ws.Range("I125").Value = 560
$ws.Range("J125").Value = 955
$ws.Range("K125").Value = 5040
$ws.Range("L125").Value = 8595
$ws.Range("M125").Value = -2580
$ws.Range("N125").Value = -13515

$ws.Range("H132").Value = 3688.5
$ws.Range("I132").Value = 3688.5
$ws.Range("K132").Value = 11065.5
$ws.Range("M132").Value = -8535.5

$ws.Range("H138").Value = 3125.9
$ws.Range("I138").Value = 3274.5454
$ws.Range("K138").Value = 9823.636200000001
$ws.Range("M138").Value = -4683.636200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6650
$ws.Range("I3").Value = 1864.6666
$ws.Range("K3").Value = 1864.6666
$ws.Range("M3").Value = -1749.6666

$ws.Range("H122").Value = 2726.2
$ws.Range("I122").Value = 2726.2
$ws.Range("K122").Value = 8178.599999999999
$ws.Range("M122").Value = -5728.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2508.9473
$ws.Range("I20").Value = 1890
$ws.Range("K20").Value = 1890
$ws.Range("M20").Value = -1643

$ws.Range("H81").Value = 4111.6665
$ws.Range("J81").Value = 4111.6665
$ws.Range("L81").Value = 4111.6665
$ws.Range("N81").Value = -6233.6665

$ws.Range("H84").Value = 4111.6665
$ws.Range("J84").Value = 4111.6665
$ws.Range("L84").Value = 12334.9995
$ws.Range("N84").Value = -22942.9995

$ws.Range("H105").Value = 3001.6086
$ws.Range("I105").Value = 2790.842
$ws.Range("K105").Value = 2790.842
$ws.Range("M105").Value = -1043.842

$ws.Range("H128").Value = 4444
$ws.Range("I128").Value = 4444
$ws.Range("K128").Value = 13332
$ws.Range("M128").Value = -10842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2099.25
$ws.Range("I16").Value = 2004.5
$ws.Range("K16").Value = 2004.5
$ws.Range("M16").Value = -1717.5

$ws.Range("H31").Value = 2215.7693
$ws.Range("I31").Value = 1338.5264
$ws.Range("K31").Value = 1338.5264
$ws.Range("M31").Value = -1043.5264

$ws.Range("H34").Value = 2215.7693
$ws.Range("I34").Value = 1338.5264
$ws.Range("K34").Value = 1338.5264
$ws.Range("M34").Value = -1136.5264

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H113").Value = 2099.25
$ws.Range("I113").Value = 2004.5
$ws.Range("K113").Value = 2004.5
$ws.Range("M113").Value = 165.5

$ws.Range("H122").Value = 2433.3333

$ws.Range("H132").Value = 4889.1816
$ws.Range("I132").Value = 4689.2
$ws.Range("K132").Value = 14067.6
$ws.Range("M132").Value = -11537.6

$ws.Range("H134").Value = 25005248
$ws.Range("I134").Value = 3000
$ws.Range("K134").Value = 9000
$ws.Range("M134").Value = -6465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 760.4
$ws.Range("I8").Value = 760.4
$ws.Range("K8").Value = 2281.2
$ws.Range("M8").Value = -2142.2

$ws.Range("H58").Value = 1987.5
$ws.Range("J58").Value = 1987.5
$ws.Range("L58").Value = 5962.5
$ws.Range("N58").Value = -6218.5

$ws.Range("H80").Value = 2779.6667
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16872

$ws.Range("H82").Value = 460668.66
$ws.Range("I82").Value = 460668.66
$ws.Range("K82").Value = 1382005.98
$ws.Range("M82").Value = -1381599.98

$ws.Range("H83").Value = 2779.6667
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54360

$ws.Range("H85").Value = 460668.66
$ws.Range("I85").Value = 460668.66
$ws.Range("K85").Value = 1382005.98
$ws.Range("M85").Value = -1380601.98

$ws.Range("H112").Value = 3802.125
$ws.Range("J112").Value = 6272.5
$ws.Range("L112").Value = 18817.5
$ws.Range("N112").Value = -21033.5

$ws.Range("H123").Value = 4600.2
$ws.Range("I123").Value = 2000.25
$ws.Range("J123").Value = 15000
$ws.Range("K123").Value = 6000.75
$ws.Range("L123").Value = 45000
$ws.Range("M123").Value = -3550.75
$ws.Range("N123").Value = -49900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6874.2383
$ws.Range("I70").Value = 7034.8667
$ws.Range("J70").Value = 6472.6665
$ws.Range("K70").Value = 7034.8667
$ws.Range("L70").Value = 6472.6665
$ws.Range("M70").Value = -6764.8667
$ws.Range("N70").Value = -7012.6665

$ws.Range("H73").Value = 6874.2383
$ws.Range("I73").Value = 7034.8667
$ws.Range("J73").Value = 6472.6665
$ws.Range("K73").Value = 7034.8667
$ws.Range("L73").Value = 6472.6665
$ws.Range("M73").Value = -6098.8667
$ws.Range("N73").Value = -8344.666499999999

$ws.Range("H133").Value = 48389
$ws.Range("J133").Value = 48389
$ws.Range("L133").Value = 48389
$ws.Range("N133").Value = -58509

$ws.Range("H138").Value = 59200
$ws.Range("J138").Value = 59200
$ws.Range("L138").Value = 59200
$ws.Range("N138").Value = -69480

$ws.Range("H139").Value = 82156.5
$ws.Range("J139").Value = 74313
$ws.Range("L139").Value = 74313
$ws.Range("N139").Value = -84593

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I25").Value = 20000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 20000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -19770
$ws.Range("N25").ClearContents()

$ws.Range("H40").Value = 2707.9
$ws.Range("I40").Value = 2012.1428
$ws.Range("K40").Value = 2012.1428
$ws.Range("M40").Value = -1876.1428

$ws.Range("H132").Value = 4651.75
$ws.Range("I132").Value = 4235.5
$ws.Range("K132").Value = 12706.5
$ws.Range("M132").Value = -10176.5

$ws.Range("H136").Value = 111118090
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 125007540
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 375022620
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -375027720

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1761
$ws.Range("I100").Value = 1038.7
$ws.Range("K100").Value = 2077.4
$ws.Range("M100").Value = -1536.4

$ws.Range("H122").Value = 2616.7896
$ws.Range("I122").Value = 4081.3333
$ws.Range("J122").Value = 1940.8462
$ws.Range("K122").Value = 12243.9999
$ws.Range("L122").Value = 5822.5386
$ws.Range("M122").Value = -9793.999899999999
$ws.Range("N122").Value = -10722.5386

$ws.Range("H126").Value = 3219.3
$ws.Range("I126").Value = 2470.5715
$ws.Range("K126").Value = 7411.7145
$ws.Range("M126").Value = -4941.7145
